$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sensor readings appended below the existing header row (AtualizaSensores fix).
$data = @(
    @("MobilityCompany", 3.25, "D1", 45262.712759699076),
    @("MobilityCompany", 3.25, "D1", 45262.713720023145),
    @("MobilityCompany", 3.25, "D1", 45262.71517917824),
    @("MobilityCompany", 3.25, "D1", 45262.71607657408),
    @("MobilityCompany", 3.25, "D1", 45262.71751649305),
    @("MobilityCompany", 3.25, "D1", 45262.718644444445)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row = $row + 1
}
